# running_time.xlsx - "Commit 3" block added (mirrors "Commit 2" layout/formulas),
# with one changed input (Other/Access count) that ripples through its formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clone the formatting of the existing "Commit 2" block (rows 34-49) onto
#    the new block (rows 52-67) so every style index (bold headers, borders,
#    totals row, etc.) matches exactly what Excel already has defined.
# ---------------------------------------------------------------------------
$ws.Range("A34:F49").Copy()
$ws.Range("A52").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Fill in the values / formulas for the new "Commit 3" section.
# ---------------------------------------------------------------------------

# Section label
$ws.Range("A52").Value = "Commit 3"

# "MARS Tool Output" header band
$ws.Range("A53").Value = "MARS Tool Output"
$ws.Range("D53").Value = "Calulations"

# "Instruction Statistics Tool" sub-header
$ws.Range("A55").Value = "Instruction Statistics Tool"

# Column headers
$ws.Range("A56").Value = "Instruction type"
$ws.Range("B56").Value = "Count"
$ws.Range("D56").Value = "Adjusted count"
$ws.Range("E56").Value = "CPI"
$ws.Range("F56").Value = "Total cycles"

# ALU
$ws.Range("A57").Value = "ALU"
$ws.Range("B57").Value = 3629
$ws.Range("D57").Formula = "=B57"
$ws.Range("E57").Value = 1
$ws.Range("F57").Formula = "=D57*E57"

# Jump
$ws.Range("A58").Value = "Jump"
$ws.Range("B58").Value = 298
$ws.Range("D58").Formula = "=B58"
$ws.Range("E58").Value = 1
$ws.Range("F58").Formula = "=D58*E58"

# Branch
$ws.Range("A59").Value = "Branch"
$ws.Range("B59").Value = 967
$ws.Range("D59").Formula = "=B59"
$ws.Range("E59").Value = 2
$ws.Range("F59").Formula = "=D59*E59"

# Memory
$ws.Range("A60").Value = "Memory"
$ws.Range("B60").Value = 621

# Other
$ws.Range("A61").Value = "Other"
$ws.Range("B61").Value = 762
$ws.Range("D61").Formula = "=B61-(B65+B66-B60)"
$ws.Range("E61").Value = 5
$ws.Range("F61").Formula = "=D61*E61"

# "Data Cache Simulation Tool" sub-header
$ws.Range("A63").Value = "Data Cache Simulation Tool"

# Column headers
$ws.Range("A64").Value = "Access"
$ws.Range("B64").Value = "Count"

# Cache hit
$ws.Range("A65").Value = "Cache hit"
$ws.Range("B65").Value = 362
$ws.Range("D65").Formula = "=B65"
$ws.Range("E65").Value = 2
$ws.Range("F65").Formula = "=D65*E65"

# Cache miss
$ws.Range("A66").Value = "Cache miss"
$ws.Range("B66").Value = 351
$ws.Range("D66").Formula = "=B66"
$ws.Range("E66").Value = 40
$ws.Range("F66").Formula = "=D66*E66"

# Total cycles for this run
$ws.Range("F67").Formula = "=SUM(F57:F66)"

# ---------------------------------------------------------------------------
# 3) Row heights for the bold "title" rows (their height nudges from 18.75 to
#    18.5 throughout the sheet, including the two new title rows).
# ---------------------------------------------------------------------------
foreach ($r in @(2, 16, 18, 32, 35, 49, 53, 67)) {
    $ws.Rows.Item($r).RowHeight = 18.5
}

# ---------------------------------------------------------------------------
# 4) View state: selection moved to B62, viewport scrolled down to row 38.
# ---------------------------------------------------------------------------
$ws.Range("B62").Select()
$excel.Goto($ws.Range("A38"), $true)
$ws.Range("B62").Select()
